# Finished skeleton of week 1 DQ1
# Fill in the "Actual time length to complete" values for the first two
# DQ1 response rows and move the active selection down to C6, matching
# the authoring session captured in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# Row 4 ("DQ1 response 1"): actual time = 0:20 -> 20/1440 of a day
$ws.Range("C4").Value = 0.013888888888888888

# Row 5 ("DQ1 response 2"): actual time = 0:06 -> 6/1440 of a day
$ws.Range("C5").Value = 0.0041666666666666666

# The Total row (C20 = SUM(C2:C19)) recalculates automatically.

# Leave the active selection on C6, where editing continued next.
$ws.Range("C6").Select()
